# Product Master issue Fix PE search
# Append 10 more "null" patient rows (154-163) to column A, mirroring the
# existing pattern of rows already present in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 154; $row -le 163; $row++) {
    $ws.Cells.Item($row, 1).Value = "null"
}
